$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.767.43'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.33%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.535.76'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.60%  '

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '616.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.27%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.70'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.96%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.529.64'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.61%  '

# Row 8
$ws.Range("E8").Value = '  +0.11%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.486'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.96%  '

# Row 10
$ws.Range("E10").Value = '  -2.18%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.91'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.55%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.431'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.39%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000223'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.30%  '

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.133.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.59%  '

# Row 15
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '32.13'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.24%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.543.93'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.52%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.761.50'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.34%  '

# Row 18
$ws.Range("E18").Value = '  +0.33%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.36'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.50%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.69%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '452.53'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.59%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.42'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.85%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.640'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.03%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.76%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.680.58'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.51%  '

# Row 27
$ws.Range("E27").Value = '  -3.65%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.50'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.14%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.37'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.18%  '

# Row 30
$ws.Range("E30").Value = '  -1.37%  '

# Row 31
$ws.Range("E31").Value = '  +0.91%  '

# Row 32
$ws.Range("E32").Value = '  -0.02%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '25.96'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.91%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.90'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.87%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.22'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.79%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.157'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.23%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.536.44'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.33%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.97'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.74%  '

# Row 39
$ws.Range("E39").Value = '  +0.02%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.02%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '176.26'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.49%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.62'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.35%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0877'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.58%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.09'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.46%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.887'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.25%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.30'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.07%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '45.83'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.55%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.58'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.70%  '

# Row 49
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.66'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.27%  '

# Row 50
$ws.Range("B50").Value = 'ONDO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.22'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.39%  '

# Row 51
$ws.Range("E51").Value = '  -2.95%  '
